# Add test case for #268: a new row 16 with a formula in A16 that
# returns "Düsseldorf" when B16 is blank (otherwise echoes B16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'

# Excel leaves the new cell selected after entry.
$ws.Range("A16").Select()
